# Insert a new row at row 10 (pushes existing row 10 and below down by one),
# then populate the new row 10 with the "capacity_to_activity" parameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above current row 10.
$ws.Rows("10:10").Insert()

# Fill in the new row with the capacity_to_activity parameter.
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_chp_biogas"
$ws.Range("C10").Value = "capacity_to_activity"
$ws.Range("D10").Value = "constant"
$ws.Range("G10").Value = 0.001
$ws.Range("H10").Value = "GW/TWh"

# Move the selection to the newly inserted row (matches the author's cursor
# position after the edit).
$ws.Range("A10:XFD10").Select()
